$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B5").Value = "1.5H"
$ws.Range("C5").Value = "0.5h"
$ws.Range("D5").Value = "0.5h"
$ws.Range("E5").Value = "Fibonacci"
$ws.Range("E5").Select()
